$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 8 (shifts old rows 8-33 down to 9-34)
$ws.Rows.Item(8).Insert()

# Insert a second brand-new row at position 28 (shifts rows 28-34 down to 29-35)
$ws.Rows.Item(28).Insert()

# --- Populate new row 8: Arcos 1855 (Cambio, nodo propio) ---
$ws.Range("A8").Value = "'-405"
$ws.Range("B8").Value = "'5/8/2025"
$ws.Range("C8").Value = "Arcos 1855"
$ws.Range("D8").Value = "'13"
$ws.Range("E8").Value = "'805791908"
$ws.Range("F8").Value = "AYKO"
$ws.Range("G8").Value = "Pendiente de Traspaso PROPIO"
$ws.Range("H8").Value = "Cambiar columna 114 picada en base, posee nodo propio.<br>"
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = "Cambio"
$ws.Range("K8").Value = "Nodo Teco"
$ws.Range("L8").Value = "Pasante"
$ws.Range("M8").Value = -58.451835
$ws.Range("N8").Value = -34.562646
$ws.Range("O8").Value = "Colegiales"
$ws.Range("P8").Value = "Capital Norte"
$ws.Range("Q8").Value = "BLO-N"
$ws.Range("R8").Value = "Fuera de Poligono OVL"

# --- Populate new row 28: CIUDAD DE LA PAZ /ALT/ 2261 ---
$ws.Range("A28").Value = "'-143"
$ws.Range("B28").Value = "'10/7/2024"
$ws.Range("C28").Value = "CIUDAD DE LA PAZ /ALT/ 2261"
$ws.Range("D28").Value = "'13"
$ws.Range("E28").Value = "'797752816"
$ws.Range("F28").Value = "AYKO"
$ws.Range("G28").Value = "Pendiente de Traspaso PROPIO"
# H28 intentionally left blank
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = "Cambio"
$ws.Range("K28").Value = "Nodo Teco"
$ws.Range("L28").Value = "Pasante"
$ws.Range("M28").Value = -58.458864
$ws.Range("N28").Value = -34.561167
$ws.Range("O28").Value = "Saavedra"
$ws.Range("P28").Value = "Capital Norte"
$ws.Range("Q28").Value = "COG-H"
$ws.Range("R28").Value = "Fuera de Poligono OVL"
